$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 57.75
$ws.Range("I9").Value = 40.5
$ws.Range("K9").Value = 40.5
$ws.Range("M9").Value = 128.5
$ws.Range("H32").Value = 1340
$ws.Range("J32").Value = 1425
$ws.Range("L32").Value = 1425
$ws.Range("N32").Value = -2077
$ws.Range("H33").Value = 137.33333
$ws.Range("I33").Value = 136.9
$ws.Range("J33").Value = 139.5
$ws.Range("K33").Value = 136.9
$ws.Range("L33").Value = 139.5
$ws.Range("M33").Value = 92.09999999999999
$ws.Range("N33").Value = -597.5
$ws.Range("H43").Value = 1000
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -931
$ws.Range("N43").ClearContents()
$ws.Range("H98").Value = 400965.03
$ws.Range("I98").Value = 622094.4
$ws.Range("K98").Value = 622094.4
$ws.Range("M98").Value = -620596.4
$ws.Range("H106").Value = 3706631.5
$ws.Range("I106").Value = 4118086.8
$ws.Range("J106").Value = 3533.3333
$ws.Range("K106").Value = 4118086.8
$ws.Range("L106").Value = 3533.3333
$ws.Range("M106").Value = -4117455.8
$ws.Range("N106").Value = -4795.3333
$ws.Range("H111").Value = 1938.174
$ws.Range("I111").Value = 1811.6471
$ws.Range("J111").Value = 2296.6667
$ws.Range("K111").Value = 5434.9413
$ws.Range("L111").Value = 6890.000100000001
$ws.Range("M111").Value = -2367.9413
$ws.Range("N111").Value = -13024.0001
$ws.Range("H113").Value = 80241.53999999999
$ws.Range("I113").Value = 103264
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 103264
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -100010
$ws.Range("N113").Value = -10008
$ws.Range("H116").Value = 9889899
$ws.Range("I116").Value = 13844849
$ws.Range("J116").Value = 2522.5
$ws.Range("K116").Value = 13844849
$ws.Range("L116").Value = 2522.5
$ws.Range("M116").Value = -13841407
$ws.Range("N116").Value = -9406.5
$ws.Range("H122").Value = 400965.03
$ws.Range("I122").Value = 622094.4
$ws.Range("K122").Value = 1866283.2
$ws.Range("M122").Value = -1863833.2
$ws.Range("H132").Value = 390052.62
$ws.Range("I132").Value = 420057.88
$ws.Range("J132").Value = 100002
$ws.Range("K132").Value = 1260173.64
$ws.Range("L132").Value = 300006
$ws.Range("M132").Value = -1257643.64
$ws.Range("N132").Value = -305066
$ws.Range("H138").Value = 1768.37
$ws.Range("I138").Value = 599.21155
$ws.Range("J138").Value = 3034.9583
$ws.Range("K138").Value = 1797.63465
$ws.Range("L138").Value = 9104.874899999999
$ws.Range("M138").Value = 3342.36535
$ws.Range("N138").Value = -19384.8749

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20540.55
$ws.Range("I32").Value = 3068.0386
$ws.Range("K32").Value = 3068.0386
$ws.Range("M32").Value = -2781.0386
$ws.Range("H36").Value = 4675.3335
$ws.Range("I36").Value = 4675.3335
$ws.Range("K36").Value = 4675.3335
$ws.Range("M36").Value = -4329.3335
$ws.Range("H45").Value = 1009.5
$ws.Range("I45").Value = 1011.5714
$ws.Range("J45").Value = 1004.6667
$ws.Range("K45").Value = 1011.5714
$ws.Range("L45").Value = 1004.6667
$ws.Range("M45").Value = -634.5714
$ws.Range("N45").Value = -1758.6667
$ws.Range("H61").Value = 2222.925
$ws.Range("I61").Value = 1519.1666
$ws.Range("K61").Value = 1519.1666
$ws.Range("M61").Value = -1307.1666
$ws.Range("H122").Value = 2352.3
$ws.Range("I122").Value = 2472.4614
$ws.Range("K122").Value = 7417.3842
$ws.Range("M122").Value = -4967.3842
$ws.Range("H132").Value = 2655.2144
$ws.Range("I132").Value = 2249.5334
$ws.Range("J132").Value = 3669.4167
$ws.Range("K132").Value = 6748.600199999999
$ws.Range("L132").Value = 11008.2501
$ws.Range("M132").Value = -4218.600199999999
$ws.Range("N132").Value = -16068.2501
$ws.Range("H136").Value = 2222.925
$ws.Range("I136").Value = 1519.1666
$ws.Range("K136").Value = 4557.4998
$ws.Range("M136").Value = -2007.4998
$ws.Range("H139").Value = 41211.5
$ws.Range("J139").Value = 45015.332
$ws.Range("L139").Value = 45015.332
$ws.Range("N139").Value = -55295.332

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5003499
$ws.Range("I7").Value = 1997.5
$ws.Range("J7").Value = 10005000
$ws.Range("K7").Value = 1997.5
$ws.Range("L7").Value = 10005000
$ws.Range("M7").Value = -1884.5
$ws.Range("N7").Value = -10005226
$ws.Range("H94").Value = 1147.25
$ws.Range("I94").Value = 871.1429000000001
$ws.Range("K94").Value = 871.1429000000001
$ws.Range("M94").Value = -420.1429000000001
$ws.Range("H99").Value = 2326
$ws.Range("I99").Value = 2157.5
$ws.Range("K99").Value = 2157.5
$ws.Range("M99").Value = -659.5
$ws.Range("H105").Value = 3190.625
$ws.Range("I105").Value = 2994.9
$ws.Range("K105").Value = 2994.9
$ws.Range("M105").Value = -1247.9
$ws.Range("H134").Value = 2475.5933
$ws.Range("I134").Value = 1608.625
$ws.Range("J134").Value = 4300.7896
$ws.Range("K134").Value = 4825.875
$ws.Range("L134").Value = 12902.3688
$ws.Range("M134").Value = -2290.875
$ws.Range("N134").Value = -17972.3688

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 542.05
$ws.Range("I22").Value = 517.4545000000001
$ws.Range("J22").Value = 572.1111
$ws.Range("K22").Value = 517.4545000000001
$ws.Range("L22").Value = 572.1111
$ws.Range("M22").Value = -167.4545000000001
$ws.Range("N22").Value = -1272.1111
$ws.Range("H31").Value = 1826.0286
$ws.Range("I31").Value = 1072.28
$ws.Range("J31").Value = 3710.4
$ws.Range("K31").Value = 1072.28
$ws.Range("L31").Value = 3710.4
$ws.Range("M31").Value = -777.28
$ws.Range("N31").Value = -4300.4
$ws.Range("H34").Value = 1826.0286
$ws.Range("I34").Value = 1072.28
$ws.Range("J34").Value = 3710.4
$ws.Range("K34").Value = 1072.28
$ws.Range("L34").Value = 3710.4
$ws.Range("M34").Value = -870.28
$ws.Range("N34").Value = -4114.4
$ws.Range("H62").Value = 20208.25
$ws.Range("I62").Value = 28162.375
$ws.Range("J62").Value = 4300
$ws.Range("K62").Value = 28162.375
$ws.Range("L62").Value = 4300
$ws.Range("M62").Value = -27538.375
$ws.Range("N62").Value = -5548
$ws.Range("H65").Value = 20208.25
$ws.Range("I65").Value = 28162.375
$ws.Range("J65").Value = 4300
$ws.Range("K65").Value = 140811.875
$ws.Range("L65").Value = 21500
$ws.Range("M65").Value = -137691.875
$ws.Range("N65").Value = -27740
$ws.Range("H68").Value = 22000
$ws.Range("J68").Value = 22000
$ws.Range("L68").Value = 22000
$ws.Range("N68").Value = -23498
$ws.Range("H71").Value = 22000
$ws.Range("J71").Value = 22000
$ws.Range("L71").Value = 66000
$ws.Range("N71").Value = -73488
$ws.Range("H99").Value = 20834666
$ws.Range("I99").Value = 62500000
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 62500000
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -62498502
$ws.Range("N99").Value = -4996
$ws.Range("H105").Value = 843.3684
$ws.Range("I105").Value = 779.1667
$ws.Range("K105").Value = 779.1667
$ws.Range("M105").Value = 967.8333
$ws.Range("H126").Value = 20834666
$ws.Range("I126").Value = 62500000
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 187500000
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -187497530
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 2445.3513
$ws.Range("I132").Value = 2078.2964
$ws.Range("K132").Value = 6234.889200000001
$ws.Range("M132").Value = -3704.889200000001
$ws.Range("H134").Value = 2036.9464
$ws.Range("I134").Value = 1382.2142
$ws.Range("J134").Value = 4001.1428
$ws.Range("K134").Value = 4146.642599999999
$ws.Range("L134").Value = 12003.4284
$ws.Range("M134").Value = -1611.642599999999
$ws.Range("N134").Value = -17073.4284

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 650
$ws.Range("J17").Value = 688.7778
$ws.Range("L17").Value = 2066.3334
$ws.Range("N17").Value = -2404.3334
$ws.Range("H131").Value = 2022.4
$ws.Range("I131").Value = 345
$ws.Range("J131").Value = 2811.7646
$ws.Range("K131").Value = 1035
$ws.Range("L131").Value = 8435.293799999999
$ws.Range("M131").Value = 4005
$ws.Range("N131").Value = -18515.2938

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 896.2941
$ws.Range("I97").Value = 717.46155
$ws.Range("K97").Value = 717.46155
$ws.Range("M97").Value = -221.46155
$ws.Range("H102").Value = 2042.9231
$ws.Range("I102").Value = 1866.8572
$ws.Range("K102").Value = 1866.8572
$ws.Range("M102").Value = -244.8571999999999
$ws.Range("H122").Value = 1390314.8
$ws.Range("I122").Value = 2779129.5
$ws.Range("K122").Value = 8337388.5
$ws.Range("M122").Value = -8334938.5
$ws.Range("H132").Value = 2888.8823
$ws.Range("I132").Value = 2782.0908
$ws.Range("J132").Value = 3084.6667
$ws.Range("K132").Value = 8346.2724
$ws.Range("L132").Value = 9254.000100000001
$ws.Range("M132").Value = -5816.2724
$ws.Range("N132").Value = -14314.0001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2829.087
$ws.Range("I7").Value = 1752.1111
$ws.Range("K7").Value = 1752.1111
$ws.Range("M7").Value = -1640.1111
$ws.Range("H16").Value = 1208.6666
$ws.Range("J16").Value = 548.8
$ws.Range("L16").Value = 548.8
$ws.Range("N16").Value = -888.8
$ws.Range("H40").Value = 4686.6665
$ws.Range("I40").Value = 3400
$ws.Range("J40").Value = 4884.615
$ws.Range("K40").Value = 3400
$ws.Range("L40").Value = 4884.615
$ws.Range("M40").Value = -3264
$ws.Range("N40").Value = -5156.615
$ws.Range("H126").Value = 2829.087
$ws.Range("I126").Value = 1752.1111
$ws.Range("K126").Value = 5256.3333
$ws.Range("M126").Value = -2786.3333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 5051394.5
$ws.Range("I107").Value = 6945167
$ws.Range("J107").Value = 1333.3334
$ws.Range("K107").Value = 20835501
$ws.Range("L107").Value = 4000.0002
$ws.Range("M107").Value = -20833581
$ws.Range("N107").Value = -7840.0002
